# Apply the two textual insertions described by the diff.
$d = $word.ActiveDocument

$quoteOpen = [char]0x2018   # ‘
$quoteClose = [char]0x201D  # ” (matches the typo in the source document)

# 1. After the first paragraph's sentence about bias/discrimination, add a new
#    sentence describing the motivation for the analysis.
$d.Content.Find.Execute("individuals. ", $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "individuals. We wanted to evaluate our data for age and gender bias.", 2)

# 2. In the EDA bullet paragraph, after the sentence about the 'EJ' feature,
#    add a new sentence calling out the 'BN' feature as a possible age value.
$d.Content.Find.Execute("among all the other features. Given", $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "among all the other features. We also considered " + $quoteOpen + "BN" + $quoteClose + " could possibly be age values. Given", 2)
